# Lecture 5 update and example
# Expanding lecture 5 and added an example for 2 moving plates
#
# 1. Retitle slide 14 ("Why multiple boundaries?" -> "Multiple boundaries?")
# 2. Remove the trailing three slides (20, 21, 22) - the "Stretching and
#    shrinking plates" / "Left stretch, right shrink" / "Move and stretch?"
#    slides - along with their notes pages.

$p = $ppt.ActivePresentation

# --- 1. Update the title text on slide 14 --------------------------------
$titleRange = $p.Slides.Item(14).Shapes.Item(1).TextFrame.TextRange

# Replace "Why multiple boundaries" with "Multiple boundaries", leaving the
# trailing "?" untouched so it keeps its own run (matches how PowerPoint
# splits runs when only part of a text box is edited).
$titleRange.Characters(1, 23).Text = "Multiple boundaries"

# --- 2. Delete the last three slides --------------------------------------
# Deleting from the end repeatedly removes slides 20, 21, and 22 (and their
# notes slides) without needing to track shifting indices.
while ($p.Slides.Count -gt 19) {
    $p.Slides.Item($p.Slides.Count).Delete()
}
